# Update cryptos list values (Price column D, Volume(1h) column E)
# generated per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.845.77'
$ws.Range("E2").Value = '  -1.94%  '
$ws.Range("D3").Value = '3.761.95'
$ws.Range("E3").Value = '  +2.56%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '620.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").Value = '3.758.78'
$ws.Range("E7").Value = '  +2.54%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +3.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.33'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.66%  '
$ws.Range("E12").Value = '  -1.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.42'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000260'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("D15").Value = '4.381.44'
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("D16").Value = '3.754.33'
$ws.Range("E16").Value = '  +2.20%  '
$ws.Range("D17").Value = '69.933.02'
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.124'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '509.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.731'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.78%  '
$ws.Range("E24").Value = '  +2.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000136'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +22.51%  '
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("E31").Value = '  +3.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.94'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.53%  '
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.23'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.340'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("E39").Value = '  +2.69%  '
$ws.Range("E40").Value = '  -3.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '429.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.08%  '
$ws.Range("D46").Value = '3.010.88'
$ws.Range("E46").Value = '  -3.84%  '
$ws.Range("E47").Value = '  -1.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.60'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.69'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.78%  '
